$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 478, shifting existing data (478:511) down to (481:514)
$ws.Rows("478:480").Insert()

# Row 478
$ws.Range("A478").Value = 2
$ws.Range("B478").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C478").Value = "Coquimbo"
$ws.Range("D478").Value = 44615
$ws.Range("E478").Value = 4
$ws.Range("F478").Value = "Fruta"
$ws.Range("G478").Value = 100102
$ws.Range("H478").Value = "Cítricos"
$ws.Range("I478").Value = 100102003
$ws.Range("J478").Value = "Limón"
$ws.Range("K478").Value = "Sin especificar"
$ws.Range("L478").Value = "1a amarillo"
$ws.Range("M478").Value = 750
$ws.Range("N478").Value = 16800
$ws.Range("O478").Value = 17000
$ws.Range("P478").Value = 16900
$ws.Range("Q478").Value = "`$/malla 16 kilos"
$ws.Range("R478").Value = "Provincia de Limarí"
$ws.Range("S478").Value = 1056
$ws.Range("T478").Value = 16

# Row 479
$ws.Range("A479").Value = 2
$ws.Range("B479").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C479").Value = "Coquimbo"
$ws.Range("D479").Value = 44615
$ws.Range("E479").Value = 4
$ws.Range("F479").Value = "Fruta"
$ws.Range("G479").Value = 100102
$ws.Range("H479").Value = "Cítricos"
$ws.Range("I479").Value = 100102003
$ws.Range("J479").Value = "Limón"
$ws.Range("K479").Value = "Sin especificar"
$ws.Range("L479").Value = "2a amarillo"
$ws.Range("M479").Value = 600
$ws.Range("N479").Value = 12800
$ws.Range("O479").Value = 13000
$ws.Range("P479").Value = 12900
$ws.Range("Q479").Value = "`$/malla 16 kilos"
$ws.Range("R479").Value = "Provincia de Limarí"
$ws.Range("S479").Value = 806
$ws.Range("T479").Value = 16

# Row 480
$ws.Range("A480").Value = 2
$ws.Range("B480").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C480").Value = "Coquimbo"
$ws.Range("D480").Value = 44615
$ws.Range("E480").Value = 4
$ws.Range("F480").Value = "Fruta"
$ws.Range("G480").Value = 100102
$ws.Range("H480").Value = "Cítricos"
$ws.Range("I480").Value = 100102003
$ws.Range("J480").Value = "Limón"
$ws.Range("K480").Value = "Sin especificar"
$ws.Range("L480").Value = "3a amarillo"
$ws.Range("M480").Value = 540
$ws.Range("N480").Value = 8800
$ws.Range("O480").Value = 9000
$ws.Range("P480").Value = 8900
$ws.Range("Q480").Value = "`$/malla 16 kilos"
$ws.Range("R480").Value = "Provincia de Limarí"
$ws.Range("S480").Value = 556
$ws.Range("T480").Value = 16

